$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fix the header: "nama" -> "name"
$ws.Range("A1").Value = "name"

# Move the active selection to F3 (as recorded by the saved view state)
$ws.Range("F3").Select()
